$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Append a new time-tracking entry as row 39 (table auto-extends) ---
# Reuse the date formatting (style) of the cell above so the new date cell
# matches the existing "Data" column formatting instead of Excel's default.
$ws.Range("H38").Copy()
$ws.Range("H39").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E39").Value = "LucaP"
$ws.Range("F39").Value = "GDPR"
$ws.Range("G39").Value = "Db, E-R"
$ws.Range("H39").Value = (Get-Date -Year 2019 -Month 4 -Day 3 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("I39").Value = 50

# --- Restore the view: scroll down and move the active selection ---
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G40").Select()
